$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 131058100
$ws.Range("B5").Value = 57725
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 102621
$ws.Range("F5").Value = "Sparvuggla"
$ws.Range("G5").Value = "Glaucidium passerinum"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("I5").Value = "'1"
$ws.Range("K5").Value = ""
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = "födosökande"
$ws.Range("N5").Value = ""
$ws.Range("P5").Value = "Österängsån, Dlr"
$ws.Range("Q5").Value = 509511
$ws.Range("R5").Value = 6751072
$ws.Range("S5").Value = 36
$ws.Range("T5").Value = "Dalarna"
$ws.Range("U5").Value = "Rättvik"
$ws.Range("V5").Value = "Dalarna"
$ws.Range("W5").Value = "Rättvik"
$ws.Range("Y5").Value = "'2026-02-07"
$ws.Range("Z5").Value = "09:39"
$ws.Range("AA5").Value = "'2026-02-07"
$ws.Range("AB5").Value = "09:39"
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AT5").Value = ""
$ws.Range("AW5").Value = "Per Lif"
$ws.Range("AX5").Value = "Per Lif"
$ws.Range("AY5").Value = ""
